$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '28.257.45'
$ws.Range('E2').Value = '  -2.58%  '
$ws.Range('D3').Value = '1.867.71'
$ws.Range('E3').Value = '  -2.13%  '
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = '318.79'
$ws.Range('E5').Value = '  -1.84%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = '0.4400'
$ws.Range('E7').Value = '  -4.16%  '
$ws.Range('D8').Value = '0.3692'
$ws.Range('E8').Value = '  -3.47%  '
$ws.Range('D9').Value = '0.07485'
$ws.Range('E9').Value = '  -2.99%  '
$ws.Range('D10').Value = '0.9360'
$ws.Range('E10').Value = '  -4.47%  '
$ws.Range('D11').Value = '21.38'
$ws.Range('E11').Value = '  -3.24%  '
$ws.Range('D12').Value = '1.887.54'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').Value = '6.693'
$ws.Range('E13').Value = '  -3.47%  '
$ws.Range('D14').Value = '5.460'
$ws.Range('E14').Value = '  -3.76%  '
$ws.Range('D15').Value = '0.06898'
$ws.Range('E15').Value = '  -2.01%  '
$ws.Range('D16').Value = '1.005'
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('D17').Value = '81.99'
$ws.Range('E17').Value = '  -2.19%  '
$ws.Range('D18').Value = '0.000009010'
$ws.Range('E18').Value = '  -4.84%  '
$ws.Range('D19').Value = '1.003'
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').Value = '15.88'
$ws.Range('E20').Value = '  -4.80%  '
$ws.Range('D21').Value = '28.238.76'
$ws.Range('E21').Value = '  -2.55%  '
$ws.Range('D22').Value = '5.111'
$ws.Range('E22').Value = '  -3.95%  '
$ws.Range('D23').Value = '10.83'
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('D24').Value = '2.100.88'
$ws.Range('E24').Value = '  -1.37%  '
$ws.Range('D25').Value = '2.024'
$ws.Range('E25').Value = '  -3.27%  '
$ws.Range('D26').Value = '155.04'
$ws.Range('E26').Value = '  -2.02%  '
$ws.Range('D27').Value = '18.35'
$ws.Range('E27').Value = '  -3.77%  '
$ws.Range('D28').Value = '5.308'
$ws.Range('E28').Value = '  -6.21%  '
$ws.Range('D29').Value = '113.44'
$ws.Range('E29').Value = '  -3.49%  '
$ws.Range('D30').Value = '1.718'
$ws.Range('E30').Value = '  -7.40%  '
$ws.Range('D31').Value = '0.09041'
$ws.Range('E31').Value = '  -2.62%  '
$ws.Range('D32').Value = '4.847'
$ws.Range('E32').Value = '  -4.38%  '
$ws.Range('D33').Value = '0.7918'
$ws.Range('E33').Value = '  -8.56%  '
$ws.Range('D34').Value = '1.171'
$ws.Range('E34').Value = '  -6.09%  '
$ws.Range('D35').Value = '2.935'
$ws.Range('E35').Value = '  -2.77%  '
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').Value = '1.130'
$ws.Range('E37').Value = '  -2.06%  '
$ws.Range('D38').Value = '0.05440'
$ws.Range('E38').Value = '  -5.11%  '
$ws.Range('E39').Value = '  -3.45%  '
$ws.Range('D40').Value = '2.950'
$ws.Range('E40').Value = '  +2.66%  '
$ws.Range('D41').Value = '0.5258'
$ws.Range('E41').Value = '  -4.53%  '
$ws.Range('D42').Value = '7.078'
$ws.Range('E42').Value = '  -4.43%  '
$ws.Range('D43').Value = '0.1682'
$ws.Range('E43').Value = '  -4.15%  '
$ws.Range('D44').Value = '8.718'
$ws.Range('E44').Value = '  -6.75%  '
$ws.Range('D45').Value = '0.06750'
$ws.Range('E45').Value = '  -1.22%  '
$ws.Range('D46').Value = '0.4876'
$ws.Range('E46').Value = '  -6.00%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = '107.57'
$ws.Range('E47').Value = '  -3.07%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '10.50'
$ws.Range('E48').Value = '  -6.49%  '
$ws.Range('D49').Value = '1.932'
$ws.Range('E49').Value = '  -5.75%  '
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('D51').Value = '1.670'
$ws.Range('E51').Value = '  -6.23%  '
